$wb = $excel.ActiveWorkbook

# --- Update status text on all three sheets: "Ready for handoff" -> "In Translation" ---

# Overview sheet: zh-cn / de-de status columns (E2, F2)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C2)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C2)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Resize the now-narrower status columns to fit the shorter text ---
# (ColumnWidth is snapped to the workbook's pixel grid, same as Excel's own
# column-width storage; 12.5 is the input that lands closest to the
# target ~13.41-character width after that snap.)

# Overview: columns E (zh-cn) and F (de-de)
$overview.Range("E1:F1").ColumnWidth = 12.5

# zh-cn / de-de: column C (Status)
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
